$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width (in "characters") of column M so the newly inserted
# column can inherit the same formatting/width, like Excel does when you
# insert a column (format copied from the column to the left).
$mWidth = $ws.Range("M1").ColumnWidth

# Insert a new blank column before column N ("Late"); this shifts the old
# N, O, P columns (Late, Outstanding label/heading, Outstanding) one to the
# right, becoming O, P, Q.
$ws.Range("N1").EntireColumn.Insert()

# Give the freshly inserted column N the same width as column M.
$ws.Range("N1").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with the new selection.
$ws.Activate() | Out-Null
$ws.Range("I17").Select() | Out-Null
